# Update countries & provincias Spain
# Applies the 12-May-2020 14:35 data refresh to the "Pais" sheet:
#  - bumps the "last updated" timestamp in A1
#  - refreshes case counters for several countries
#  - Catar overtakes Singapur/Bielorrusia in the ranking (rows 28-30)
#  - Cabo Verde overtakes Zambia/Etiopia in the ranking (rows 139-141)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 14:35"

# --- Helper: write Total/Nuevos/Activos/Recuperados/Criticos/MuertesHoy/Muertes
function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Espana (row 5)
Set-Row 5 269520 1377 180470 62130 1534 176 26920

# India (row 15)
Set-Row 15 71441 673 23049 46082 0 16 2310

# Paises Bajos (row 19)
Set-Row 19 42984 196 0 37224 498 54 5510

# Suiza (row 23)
Set-Row 23 30380 36 26800 1723 89 12 1857

# Portugal (row 26)
Set-Row 26 27913 234 3013 23737 113 19 1163

# Suecia (row 27)
Set-Row 27 27272 602 4971 18988 360 57 3313

# Rows 28-30 re-rank: Catar passes Singapur and Bielorrusia.
$ws.Range("A28").Value = "Catar"
Set-Row 28 25149 1526 3019 22116 72 0 14

$ws.Range("A29").Value = "Singapur"
Set-Row 29 24671 884 3225 21425 24 0 21

$ws.Range("A30").Value = "Bielorrusia"
Set-Row 30 23906 0 6531 17240 92 0 135

# Finlandia (row 57)
Set-Row 57 6003 19 4300 1428 34 4 275

# Croacia (row 77)
Set-Row 77 2207 11 1808 308 11 0 91

# Estado de Palestina (row 129)
Set-Row 129 375 0 316 57 0 0 2

# Rows 139-141 re-rank: Cabo Verde passes Zambia and Etiopia.
$ws.Range("A139").Value = "Cabo Verde"
Set-Row 139 267 7 58 207 0 0 2

$ws.Range("A140").Value = "Zambia"
Set-Row 140 267 0 117 143 1 0 7

$ws.Range("A141").Value = "Etiopia"
Set-Row 141 261 11 106 150 0 0 5
